# Add api endpoint validation
# Appends one new data row (row 54) to each of the four worksheets,
# mirroring the existing row 53 pattern but with the hour advanced by one.

$wb = $excel.ActiveWorkbook

function Add-SensorRow {
    param($ws, $timeVal, $colB, $colC, $colD, $colE, $colF, $colG, $colH, $colI)

    $row = 54

    $ws.Cells.Item($row, 1).Value = $timeVal
    $ws.Cells.Item($row, 2).Value = $colB
    $ws.Cells.Item($row, 3).Value = $colC
    $ws.Cells.Item($row, 4).Value = $colD
    $ws.Cells.Item($row, 5).Value = $colE
    $ws.Cells.Item($row, 6).Value = $colF

    # Column G holds a 24-digit id string; forcing text format keeps the
    # exact digits instead of Excel rounding it to a double, and resetting
    # the style back to Normal afterwards avoids leaving a stray custom
    # number format behind on the cell.
    $ws.Cells.Item($row, 7).NumberFormat = "@"
    $ws.Cells.Item($row, 7).Value = $colG
    $ws.Cells.Item($row, 7).Style = "Normal"

    $ws.Cells.Item($row, 8).Value = $colH
    $ws.Cells.Item($row, 9).Value = $colI
}

# Sheet 1: ROW35-FE-LIFTER
$ws1 = $wb.Worksheets.Item("ROW35-FE-LIFTER")
Add-SensorRow $ws1 "2025-03-06 13:42:06" "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c," "0x01,0x90," "0x d" 400 "568631262647113770877196" 400 13

# Sheet 2: ROW35-MID-LIFTER
$ws2 = $wb.Worksheets.Item("ROW35-MID-LIFTER")
Add-SensorRow $ws2 "2025-03-06 13:29:35" "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c," "0x01,0x90," "0x e" 400 "568631262647113770942732" 400 14

# Sheet 3: ROW02-FE-LIFTER
$ws3 = $wb.Worksheets.Item("ROW02-FE-LIFTER")
Add-SensorRow $ws3 "2025-03-06 13:51:45" "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c," "0x01,0x90," "0xff" 400 "568631262647113769959692" 400 255

# Sheet 4: ROW02-MID-LIFTER
$ws4 = $wb.Worksheets.Item("ROW02-MID-LIFTER")
Add-SensorRow $ws4 "2025-03-06 13:41:15" "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c," "0x01,0x90," "0x 3" 400 "568631262647113769959692" 400 3
